$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell style (bold, centered, bordered) from AC1 onto the
# three new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Add the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 85   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 77   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
